# Snow timeline update: add Canada sites (queens, old_jack_pine) timeseries rows
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, location, date serial, FSC value, melt year
$newRows = @(
    @(1279, "queens", 44634, 98, 2022),
    @(1280, "queens", 44635, 96, 2022),
    @(1281, "queens", 44636, 95, 2022),
    @(1282, "queens", 44637, 94, 2022),
    @(1283, "queens", 44638, 93, 2022),
    @(1284, "queens", 44639, 91, 2022),
    @(1285, "queens", 44640, 90, 2022),
    @(1286, "queens", 44641, 87, 2022),
    @(1287, "queens", 44642, 83, 2022),
    @(1288, "queens", 44643, 79, 2022),
    @(1289, "queens", 44644, 50, 2022),
    @(1290, "queens", 44645, 25, 2022),
    @(1291, "queens", 44646, 12, 2022),
    @(1292, "queens", 44647, 7, 2022),
    @(1293, "queens", 44648, 6, 2022),
    @(1294, "queens", 44649, 5, 2022),
    @(1295, "queens", 44265, 97, 2021),
    @(1296, "queens", 44266, 90, 2021),
    @(1297, "queens", 44267, 75, 2021),
    @(1298, "queens", 44268, 45, 2021),
    @(1299, "queens", 44269, 40, 2021),
    @(1300, "queens", 44270, 38, 2021),
    @(1301, "queens", 44271, 37, 2021),
    @(1302, "queens", 44272, 25, 2021),
    @(1303, "queens", 44273, 10, 2021),
    @(1304, "queens", 44274, 7, 2021),
    @(1305, "queens", 44275, 4, 2021),
    @(1306, "queens", 43547, 95, 2019),
    @(1307, "queens", 43548, 85, 2019),
    @(1308, "queens", 43549, 84, 2019),
    @(1309, "queens", 43550, 80, 2019),
    @(1310, "queens", 43551, 65, 2019),
    @(1311, "queens", 43552, 50, 2019),
    @(1312, "queens", 43553, 30, 2019),
    @(1313, "queens", 43554, 15, 2019),
    @(1314, "queens", 43556, 10, 2019),
    @(1315, "queens", 43557, 5, 2019),
    @(1316, "old_jack_pine", 44673, 100, 2022),
    @(1317, "old_jack_pine", 44674, 99, 2022),
    @(1318, "old_jack_pine", 44675, 98, 2022),
    @(1319, "old_jack_pine", 44676, 97, 2022),
    @(1320, "old_jack_pine", 44677, 96, 2022),
    @(1321, "old_jack_pine", 44678, 94, 2022),
    @(1322, "old_jack_pine", 44679, 90, 2022),
    @(1323, "old_jack_pine", 44680, 80, 2022),
    @(1324, "old_jack_pine", 44681, 65, 2022),
    @(1325, "old_jack_pine", 44682, 60, 2022),
    @(1326, "old_jack_pine", 44683, 50, 2022),
    @(1327, "old_jack_pine", 44684, 50, 2022),
    @(1328, "old_jack_pine", 44685, 35, 2022),
    @(1329, "old_jack_pine", 44686, 20, 2022),
    @(1330, "old_jack_pine", 44687, 8, 2022),
    @(1331, "old_jack_pine", 44688, 5, 2022),
    @(1332, "old_jack_pine", 44287, 100, 2021),
    @(1333, "old_jack_pine", 44288, 99, 2021),
    @(1334, "old_jack_pine", 44289, 98, 2021),
    @(1335, "old_jack_pine", 44290, 96, 2021),
    @(1336, "old_jack_pine", 44291, 94, 2021),
    @(1337, "old_jack_pine", 44292, 90, 2021),
    @(1338, "old_jack_pine", 44293, 75, 2021),
    @(1339, "old_jack_pine", 44294, 60, 2021),
    @(1340, "old_jack_pine", 44295, 54, 2021),
    @(1341, "old_jack_pine", 44296, 45, 2021),
    @(1342, "old_jack_pine", 44299, 44, 2021),
    @(1343, "old_jack_pine", 44300, 35, 2021),
    @(1344, "old_jack_pine", 44301, 20, 2021),
    @(1345, "old_jack_pine", 44302, 15, 2021),
    @(1346, "old_jack_pine", 44303, 10, 2021),
    @(1347, "old_jack_pine", 43944, 95, 2020),
    @(1348, "old_jack_pine", 43945, 90, 2020),
    @(1349, "old_jack_pine", 43946, 75, 2020),
    @(1350, "old_jack_pine", 43947, 35, 2020),
    @(1351, "old_jack_pine", 43948, 15, 2020),
    @(1352, "old_jack_pine", 43949, 5, 2020),
    @(1353, "old_jack_pine", 43566, 97, 2019),
    @(1354, "old_jack_pine", 43567, 95, 2019),
    @(1355, "old_jack_pine", 43568, 90, 2019),
    @(1356, "old_jack_pine", 43569, 87, 2019),
    @(1357, "old_jack_pine", 43570, 85, 2019),
    @(1358, "old_jack_pine", 43571, 73, 2019),
    @(1359, "old_jack_pine", 43572, 60, 2019),
    @(1360, "old_jack_pine", 43573, 33, 2019),
    @(1361, "old_jack_pine", 43574, 20, 2019),
    @(1362, "old_jack_pine", 43575, 10, 2019),
    @(1363, "old_jack_pine", 43576, 7, 2019),
    @(1364, "old_jack_pine", 43577, 4, 2019)
)

foreach ($r in $newRows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Match the author's final cursor position/selection after pasting the new data
$ws.Range("C1365").Select() | Out-Null

